# Rename the inline picture shapes ("image2.png" <-> "image1.png" /
# "image1.jpg" <-> "image2.jpg") that live in the document's footers and
# header, per the commit's docPr/cNvPr @name swap.
#
# InlineShape objects don't expose the drawing's name directly on the
# header/footer story Range (InlineShapes(1) against the whole story
# fails to address the picture), so each picture is located through the
# specific paragraph that contains it before its .Name is set.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Set-PictureName($headerFooter, $newName) {
    $paraCount = $headerFooter.Range.Paragraphs.Count
    for ($i = 1; $i -le $paraCount; $i++) {
        $para = $headerFooter.Range.Paragraphs($i)
        if ($para.Range.InlineShapes.Count -gt 0) {
            $shape = $para.Range.InlineShapes(1)
            $shape.Name = $newName
        }
    }
}

# Footer 1 (first footer) and Footer 2 (second footer) both hold the
# Pearson Edexcel logo: image2.png -> image1.png
Set-PictureName $sec.Footers(1) "image1.png"
Set-PictureName $sec.Footers(2) "image1.png"

# Header 2 holds the BTEC logo: image1.jpg -> image2.jpg
Set-PictureName $sec.Headers(2) "image2.jpg"
